# Insert a new weekly price record at row 35 (pushing the existing rows
# 35-140 down to 36-141), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 35..140 down by one, creating space for the new record.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly data point.
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44690
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112012
$ws.Range("G35").Value = "Espinaca"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 70
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 9000
$ws.Range("M35").Value = 9000
$ws.Range("N35").Value = "$/docena de atados"
$ws.Range("O35").Value = "Región de La Araucanía"
$ws.Range("P35").Value = 3000
$ws.Range("Q35").Value = 3
$ws.Range("R35").Value = "Hortaliza"
